$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    # Force the cell to keep a numeric-looking string as TEXT (matches the
    # source data, which stores prices/volumes as inline strings, not
    # numbers) instead of letting Excel auto-convert it to a number.
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "66.272.78"
Set-TextValue "E2" "  -1.37%  "

# Row 3 - Ethereum
Set-TextValue "D3" "3.563.64"
Set-TextValue "E3" "  +1.46%  "

# Row 4 - TetherUSD
Set-TextValue "D4" "0.999"
Set-TextValue "E4" "  -0.02%  "

# Row 5 - BNB
Set-TextValue "D5" "608.09"
Set-TextValue "E5" "  -0.41%  "

# Row 6 - Solana
Set-TextValue "D6" "144.64"
Set-TextValue "E6" "  -2.33%  "

# Row 7 - LidoStakedEther
Set-TextValue "D7" "3.561.81"
Set-TextValue "E7" "  +1.45%  "

# Row 8 - USDC
Set-TextValue "E8" "  -0.03%  "

# Row 9 - XRP
Set-TextValue "D9" "0.483"
Set-TextValue "E9" "  +0.62%  "

# Row 10 - was Toncoin, now Dogecoin
$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
Set-TextValue "D10" "0.137"
Set-TextValue "E10" "  -3.80%  "

# Row 11 - was Dogecoin, now Toncoin
$ws.Range("B11").Value = "Toncoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue "D11" "8.08"
Set-TextValue "E11" "  +0.30%  "

# Row 12 - Cardano
Set-TextValue "D12" "0.412"
Set-TextValue "E12" "  -2.57%  "

# Row 13 - WrappedliquidstakedEther2.0
Set-TextValue "D13" "4.167.12"

# Row 14 - ShibaInu
Set-TextValue "E14" "  -3.93%  "

# Row 15 - Avalanche
Set-TextValue "D15" "30.21"
Set-TextValue "E15" "  -4.23%  "

# Row 16 - WrappedEther
Set-TextValue "D16" "3.557.73"
Set-TextValue "E16" "  +1.36%  "

# Row 17 - WrappedBTC
Set-TextValue "D17" "66.354.74"
Set-TextValue "E17" "  -1.33%  "

# Row 18 - TRON
Set-TextValue "E18" "  -1.09%  "

# Row 19 - Uniswap
Set-TextValue "D19" "11.36"
Set-TextValue "E19" "  +4.18%  "

# Row 20 - Polkadot
Set-TextValue "D20" "6.22"
Set-TextValue "E20" "  -2.25%  "

# Row 21 - Chainlink
Set-TextValue "D21" "14.95"
Set-TextValue "E21" "  -3.12%  "

# Row 22 - BitcoinCash
Set-TextValue "D22" "429.68"
Set-TextValue "E22" "  -1.51%  "

# Row 23 - Polygon
Set-TextValue "D23" "0.606"
Set-TextValue "E23" "  -0.65%  "

# Row 24 - Litecoin
Set-TextValue "E24" "  -1.48%  "

# Row 25 - WrappedeETH
Set-TextValue "D25" "3.700.64"
Set-TextValue "E25" "  +1.26%  "

# Row 26 - Dai
Set-TextValue "E26" "  -0.04%  "

# Row 27 - PEPE
Set-TextValue "D27" "0.0000122"
Set-TextValue "E27" "  +2.35%  "

# Row 28 - RenderToken
Set-TextValue "D28" "8.10"
Set-TextValue "E28" "  -2.27%  "

# Row 29 - InternetComputer(DFINITY)
Set-TextValue "D29" "9.23"
Set-TextValue "E29" "  -6.49%  "

# Row 30 - PancakeSwap
Set-TextValue "E30" "  -1.38%  "

# Row 31 - Binance-PegBSC-USD
Set-TextValue "E31" "  -0.04%  "

# Row 32 - Fetch.AI
Set-TextValue "E32" "  -5.83%  "

# Row 33 - Kaspa
Set-TextValue "D33" "0.158"
Set-TextValue "E33" "  -4.41%  "

# Row 34 - EthereumClassic
Set-TextValue "D34" "25.49"
Set-TextValue "E34" "  -0.39%  "

# Row 35 - RenzoRestakedETH
Set-TextValue "D35" "3.553.06"
Set-TextValue "E35" "  +1.35%  "

# Row 36 - USDe
Set-TextValue "E36" "  -0.03%  "

# Row 37 - ImmutableX
Set-TextValue "D37" "1.76"
Set-TextValue "E37" "  -3.21%  "

# Row 38 - was NEARProtocol, now Aptos
$ws.Range("B38").Value = "Aptos"
$ws.Range("C38").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue "D38" "7.83"
Set-TextValue "E38" "  -2.51%  "

# Row 39 - was Aptos, now NEARProtocol
$ws.Range("B39").Value = "NEARProtocol"
$ws.Range("C39").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D39" "5.63"
Set-TextValue "E39" "  -5.61%  "

# Row 40 - FirstDigitalUSD
Set-TextValue "D40" "0.999"
Set-TextValue "E40" "  +0.00%  "

# Row 41 - Monero
Set-TextValue "D41" "174.86"
Set-TextValue "E41" "  -0.43%  "

# Row 42 - Hedera
Set-TextValue "D42" "0.0860"
Set-TextValue "E42" "  -4.94%  "

# Row 43 - Filecoin
Set-TextValue "D43" "5.32"
Set-TextValue "E43" "  -1.64%  "

# Row 44 - Mantle
Set-TextValue "D44" "0.896"
Set-TextValue "E44" "  +0.00%  "

# Row 45 - Stacks
Set-TextValue "D45" "1.91"
Set-TextValue "E45" "  -7.00%  "

# Row 47 - ONDO
Set-TextValue "E47" "  -1.47%  "

# Row 48 - InjectiveProtocol
Set-TextValue "D48" "26.04"
Set-TextValue "E48" "  -9.12%  "

# Row 49 - dogwifhat
Set-TextValue "E49" "  -2.57%  "

# Row 50 - Cosmos
Set-TextValue "D50" "7.15"
Set-TextValue "E50" "  -4.35%  "

# Row 51 - EnergySwap
Set-TextValue "D51" "23.02"
Set-TextValue "E51" "  +6.10%  "
